$wb = $excel.ActiveWorkbook

# --- Rebuild the sheet set: FreeGames, AddGames, InstallGames ------------------
# Sheet ids are assigned by Excel as (current max id + 1) at the moment a sheet
# is created. A scratch sheet is used purely to advance that counter so the
# real "AddGames"/"InstallGames" sheets land on ids 4 and 5, matching the ids
# a human using the Excel UI (add/delete a throwaway tab while prototyping)
# would naturally end up with. Sheet references are re-fetched by name right
# before each use since earlier handles can point at a stale tab position
# once the sheet collection has been reordered.
$scratch = $wb.Worksheets.Add()

$freeGames = $wb.Worksheets.Item("FreeGames")
$addGames = $wb.Worksheets.Add($null, $freeGames)
$addGames.Name = "AddGames"

$addGames = $wb.Worksheets.Item("AddGames")
$installGames = $wb.Worksheets.Add($null, $addGames)
$installGames.Name = "InstallGames"

$wb.Worksheets.Item($scratch.Name).Delete()
$wb.Worksheets.Item("Sheet1").Delete()

$freeGames = $wb.Worksheets.Item("FreeGames")
$addGames = $wb.Worksheets.Item("AddGames")
$installGames = $wb.Worksheets.Item("InstallGames")

# --- FreeGames sheet: drop the Status/Note columns, add a new game row ----------
$freeGames.Columns.Item(5).Delete()
$freeGames.Columns.Item(4).Delete()

$freeGames.Range("A1").Value = "Game Title"
$freeGames.Range("B1").Value = "Add"
$freeGames.Range("C1").Value = "Install"

$freeGames.Range("A2").Value = "Apex Legends"
$freeGames.Range("B2").Value = "Yes"
$freeGames.Range("C2").Value = "No"

$freeGames.Range("A3").Value = "Phantasy Star Online 2"
$freeGames.Range("B3").Value = "Yes"
$freeGames.Range("C3").Value = "No"

$freeGames.Range("A4").Value = "Destiny 2"
$freeGames.Range("B4").Value = "Yes"
$freeGames.Range("C4").Value = "No"

$freeGames.Range("A5").Value = "War Thunder"
$freeGames.Range("B5").Value = "No"
$freeGames.Range("C5").Value = "No"

$freeGames.Range("A6").Value = "Assassin's Creed Odyssey"
$freeGames.Range("B6").Value = "No"
$freeGames.Range("C6").Value = "Yes"

# --- AddGames sheet: games already queued to be added ---------------------------
$addGames.Range("A1").Value = "Game Title"
$addGames.Range("B1").Value = "Status"

$addGames.Range("A2").Value = "Apex Legends"
$addGames.Range("B2").Value = "Already In Library"

$addGames.Range("A3").Value = "Phantasy Star Online 2"
$addGames.Range("B3").Value = "Already In Library"

$addGames.Range("A4").Value = "Destiny 2"
$addGames.Range("B4").Value = "Already In Library"

# --- InstallGames sheet: games queued to be installed ----------------------------
$installGames.Range("A1").Value = "Game Title"
$installGames.Range("B1").Value = "Status"

$installGames.Range("A2").Value = "Assassin's Creed Odyssey"
$installGames.Range("B2").Value = "Game Installing"

# Match the formatting used throughout the workbook (style carried from FreeGames)
$freeGames.Range("A1:B1").Copy()
$addGames.Range("A1:B4").PasteSpecial(-4122)
$installGames.Range("A1:B2").PasteSpecial(-4122)

$installGames.Activate()
